# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# with latest crypto market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "69.503.21"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.691.45"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "679.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.58"
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000234"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "4.312.45"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.49"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "3.680.55"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "69.448.81"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "3.838.11"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.58"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.01"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "3.679.93"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.163"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.24"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0902"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.85%  "
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.88%  "
